# Corrected Calibration and Ingest Sheets for Coastal Gliders
# - FLORT CC_scattering_angle -> 124
# - FLORT CC_angular_resolution -> 1.076
# - Asset_Cal_Info becomes the active/selected sheet (cell F29 selected)

$wb = $excel.ActiveWorkbook

$calInfo = $wb.Worksheets.Item("Asset_Cal_Info")

# FLORT calibration coefficients (rows identified by their
# "Calibration Cofficient Name" in column E):
#   row 4 -> CC_scattering_angle    117  -> 124
#   row 6 -> CC_angular_resolution  1.08 -> 1.076
$calInfo.Range("F4").Value = 124
$calInfo.Range("F6").Value = 1.076

# Make Asset_Cal_Info the active sheet/tab and move the selection to F29,
# matching the saved workbook view state.
$calInfo.Activate()
$calInfo.Range("F29").Select()
